# Trade #51 closed at 2026-02-17 12:49:11 - unknown UNKNOWN +0.000%
#
# 1) Summary sheet: bump Total Trades / recompute Win Rate %
# 2) Strategy Status sheet: bump MarketMaking Trades / recompute Win Rate %
# 3) All Trades + MarketMaking sheets: append the new closed trade as row 52

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 51
$summary.Range("B9").Value = 41.18

# ---------------------------------------------------------------------
# Strategy Status (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 51
$status.Range("G4").Value = 41.18

# ---------------------------------------------------------------------
# Helper: write the new trade row (#51 / spreadsheet row 52) onto a
# trade-log sheet ("All Trades" and "MarketMaking" share the same
# layout/content).
# ---------------------------------------------------------------------
function Add-ClosedTrade($ws) {
    $row = 52

    $ws.Cells.Item($row, 1).Value = 51

    # Date/time columns hold plain text like "2026-02-17" - Excel's
    # auto-detection would otherwise coerce that into a date serial, so
    # force the cell to Text first and restore the default style after.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "12:49:05"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.97
    $ws.Cells.Item($row, 7).Value = 0.97
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.16
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-ClosedTrade $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-ClosedTrade $marketMaking
